$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column D (shifts old D:K -> F:M)
$ws.Columns("D:E").Insert()

# Copy formatting from the (now-shifted) old column D, which landed in F,
# onto the two newly inserted blank columns D:E so they inherit the same
# per-row number formats/styles (date header row, numeric data rows, etc.)
$ws.Range("F7:F102").Copy()
$ws.Range("D7:E102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# The copy/paste above also stamped formatting onto the blank spacer rows
# (36 and 78) that have no real content in any column - clear those back out
# so those rows stay absent, matching the original sheet's row layout.
$ws.Range("D36:E36").Clear()
$ws.Range("D78:E78").Clear()

# Populate the two new columns with the new quarter's figures
$ws.Range("D7").Value2 = 43464
$ws.Range("E7").Value2 = 43373
$ws.Range("D8").Value2 = 748400
$ws.Range("E8").Value2 = 725300
$ws.Range("D9").Value2 = 459600
$ws.Range("E9").Value2 = 446200
$ws.Range("D10").Value2 = 288800
$ws.Range("E10").Value2 = 279100
$ws.Range("D12").Value2 = "NA"
$ws.Range("E12").Value2 = "NA"
$ws.Range("D13").Value2 = 0
$ws.Range("E13").Value2 = 0
$ws.Range("D14").Value2 = "NA"
$ws.Range("E14").Value2 = "NA"
$ws.Range("D15").Value2 = 0
$ws.Range("E15").Value2 = 0
$ws.Range("D17").Value2 = 637200
$ws.Range("E17").Value2 = 619800
$ws.Range("D18").Value2 = 111200
$ws.Range("E18").Value2 = 105500
$ws.Range("D20").Value2 = -4100
$ws.Range("E20").Value2 = -5300
$ws.Range("D21").Value2 = 136400
$ws.Range("E21").Value2 = 127500
$ws.Range("D22").Value2 = 0
$ws.Range("E22").Value2 = 0
$ws.Range("D23").Value2 = 107100
$ws.Range("E23").Value2 = 100200
$ws.Range("D24").Value2 = 16800
$ws.Range("E24").Value2 = 10500
$ws.Range("D25").Value2 = 0
$ws.Range("E25").Value2 = 0
$ws.Range("D26").Value2 = 90300
$ws.Range("E26").Value2 = 89700
$ws.Range("D27").Value2 = 90300
$ws.Range("E27").Value2 = 89700
$ws.Range("D28").Value2 = 0
$ws.Range("E28").Value2 = 0
$ws.Range("D29").Value2 = 800
$ws.Range("E29").Value2 = 600
$ws.Range("D30").Value2 = 0
$ws.Range("E30").Value2 = 0
$ws.Range("D31").Value2 = 0
$ws.Range("E31").Value2 = 0
$ws.Range("D32").Value2 = 4100
$ws.Range("E32").Value2 = 5300
$ws.Range("D33").Value2 = 91100
$ws.Range("E33").Value2 = 90300
$ws.Range("D34").Value2 = 0
$ws.Range("E34").Value2 = 0
$ws.Range("D35").Value2 = 91100
$ws.Range("E35").Value2 = 90300
$ws.Range("D38").Value2 = 43464
$ws.Range("E38").Value2 = 43373
$ws.Range("D41").Value2 = 142500
$ws.Range("E41").Value2 = 126100
$ws.Range("D42").Value2 = 0
$ws.Range("E42").Value2 = 0
$ws.Range("D43").Value2 = 561800
$ws.Range("E43").Value2 = 551800
$ws.Range("D44").Value2 = 364300
$ws.Range("E44").Value2 = 377600
$ws.Range("D45").Value2 = 45800
$ws.Range("E45").Value2 = 54100
$ws.Range("D46").Value2 = 1114400
$ws.Range("E46").Value2 = 1109600
$ws.Range("D47").Value2 = 0
$ws.Range("E47").Value2 = 0
$ws.Range("D48").Value2 = 442600
$ws.Range("E48").Value2 = 446100
$ws.Range("D49").Value2 = 2079500
$ws.Range("E49").Value2 = 2113300
$ws.Range("D50").Value2 = 0
$ws.Range("E50").Value2 = 0
$ws.Range("D51").Value2 = 0
$ws.Range("E51").Value2 = 0
$ws.Range("D52").Value2 = 172800
$ws.Range("E52").Value2 = 244200
$ws.Range("D53").Value2 = 0
$ws.Range("E53").Value2 = 0
$ws.Range("D54").Value2 = 3809300
$ws.Range("E54").Value2 = 3913200
$ws.Range("D57").Value2 = 227800
$ws.Range("E57").Value2 = 205400
$ws.Range("D58").Value2 = 138300
$ws.Range("E58").Value2 = 14300
$ws.Range("D59").Value2 = 354700
$ws.Range("E59").Value2 = 368300
$ws.Range("D60").Value2 = 720800
$ws.Range("E60").Value2 = 588000
$ws.Range("D61").Value2 = 612300
$ws.Range("E61").Value2 = 823000
$ws.Range("D62").Value2 = 246500
$ws.Range("E62").Value2 = 285400
$ws.Range("D63").Value2 = 0
$ws.Range("E63").Value2 = 0
$ws.Range("D64").Value2 = 0
$ws.Range("E64").Value2 = 0
$ws.Range("D65").Value2 = 0
$ws.Range("E65").Value2 = 0
$ws.Range("D66").Value2 = 1579600
$ws.Range("E66").Value2 = 1696400
$ws.Range("D68").Value2 = 0
$ws.Range("E68").Value2 = 0
$ws.Range("D69").Value2 = 0
$ws.Range("E69").Value2 = 0
$ws.Range("D70").Value2 = 0
$ws.Range("E70").Value2 = 0
$ws.Range("D71").Value2 = 0
$ws.Range("E71").Value2 = 0
$ws.Range("D72").Value2 = 2523700
$ws.Range("E72").Value2 = 2432500
$ws.Range("D73").Value2 = 0
$ws.Range("E73").Value2 = 0
$ws.Range("D74").Value2 = 0
$ws.Range("E74").Value2 = 0
$ws.Range("D75").Value2 = 0
$ws.Range("E75").Value2 = 0
$ws.Range("D76").Value2 = 2229700
$ws.Range("E76").Value2 = 2216800
$ws.Range("D77").Value2 = 0
$ws.Range("E77").Value2 = 0
$ws.Range("D80").Value2 = 43464
$ws.Range("E80").Value2 = 43373
$ws.Range("D81").Value2 = 91100
$ws.Range("E81").Value2 = 90300
$ws.Range("D83").Value2 = 29300
$ws.Range("E83").Value2 = 27300
$ws.Range("D84").Value2 = 0
$ws.Range("E84").Value2 = 0
$ws.Range("D85").Value2 = 0
$ws.Range("E85").Value2 = 0
$ws.Range("D86").Value2 = 0
$ws.Range("E86").Value2 = 0
$ws.Range("D87").Value2 = 0
$ws.Range("E87").Value2 = 0
$ws.Range("D88").Value2 = 0
$ws.Range("E88").Value2 = 0
$ws.Range("D89").Value2 = 125500
$ws.Range("E89").Value2 = 141900
$ws.Range("D91").Value2 = -18700
$ws.Range("E91").Value2 = -20900
$ws.Range("D92").Value2 = 0
$ws.Range("E92").Value2 = 0
$ws.Range("D93").Value2 = 0
$ws.Range("E93").Value2 = 0
$ws.Range("D94").Value2 = -21700
$ws.Range("E94").Value2 = -20400
$ws.Range("D96").Value2 = 0
$ws.Range("E96").Value2 = 0
$ws.Range("D97").Value2 = 0
$ws.Range("E97").Value2 = 0
$ws.Range("D98").Value2 = 0
$ws.Range("E98").Value2 = 0
$ws.Range("D99").Value2 = 0
$ws.Range("E99").Value2 = 0
$ws.Range("D100").Value2 = -80800
$ws.Range("E100").Value2 = -99600
$ws.Range("D101").Value2 = -6600
$ws.Range("E101").Value2 = 2800
$ws.Range("D102").Value2 = 16400
$ws.Range("E102").Value2 = 24700
